# Apply updates to the "Blogs_used_list" worksheet per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = "error_occured"
$ws.Range("E4").Value = "https://www.error_link.com"

# Row 8
$ws.Range("D8").Value = "error_occured"
$ws.Range("E8").Value = "https://www.error_link.com"

# Row 28 - title text updated (remove trailing " - (미완)")
$ws.Range("D28").Value = "[keep9oing] Counterfactual Multi-Agent Policy Gradient (COMA) 리뷰"

# Row 29
$ws.Range("D29").Value = "error_occured"
$ws.Range("E29").Value = "https://www.error_link.com"

# Row 32 - title and link changed to a different article
$ws.Range("D32").Value = "DFS(Depth First Search), BFS(Breadth First Search) - 깊이/너비 우선 탐색"
$ws.Range("E32").Value = "https://dodonam.tistory.com/290"

# Row 35
$ws.Range("D35").Value = "error_occured"
$ws.Range("E35").Value = "https://www.error_link.com"

# Row 36
$ws.Range("D36").Value = "error_occured"
$ws.Range("E36").Value = "https://www.error_link.com"

# Row 37
$ws.Range("D37").Value = "error_occured"
$ws.Range("E37").Value = "https://www.error_link.com"
